# 2021-05-13 문구수정 및 vali 변경 commit
#
# The "example / validation hint" row (old row 2) gets folded into the
# header row (row 1) as a second line of text (e.g. "구분" becomes
# "구분`n신규 = 1`n경력 = 2"), the old row is removed (old row 3's sample
# data shifts up to become row 2), and the sample name is updated from
# 신동환 to 홍길동.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the old "example value" row. The real data row (old row 3)
#    shifts up to become row 2.
$ws.Rows.Item(2).Delete()

# 2) Re-word the header cells so each now carries its old validation /
#    example text as an embedded second (and further) line.
$ws.Range("A1").Value = "구분`n신규 = 1`n경력 = 2"
$ws.Range("C1").Value = "주민등록번호`n000000-0000000"
$ws.Range("D1").Value = "휴대폰번호`n000-0000-0000"
$ws.Range("E1").Value = "주소`n서울 = 1`n경기 = 2`n충청북도 = 3`n충청남도 = 4`n강원도 = 5`n경상북도 = 6`n경상남도 = 7`n전라북도 = 8`n전라남도 = 9`n인천 = 10`n세종 = 11`n대전 = 12`n대구 = 13`n울산 = 14`n광주 = 15`n부산 = 16`n제주 = 17"
$ws.Range("F1").Value = "금융상품유형`n대출 = 1`n시설대여 및 연불판매 = 2`n할부 = 3`n어음할인 = 4`n매출채권 매입 = 5`n지급보증 = 6`n기타 대출성상품 = 7"
$ws.Range("H1").Value = "법인등록번호`n000000-0000000"
$ws.Range("I1").Value = "교육이수번호/인증서번호`n0000000000"
$ws.Range("J1").Value = "경력시작일`nYYYY-MM-DD"
$ws.Range("K1").Value = "경력종료일`nYYYY-MM-DD"
$ws.Range("L1").Value = "계약일자`nYYYY-MM-DD"
$ws.Range("M1").Value = "위탁예정기간`nYYYY-MM-DD"
# B1 (성명) and G1 (법인명) keep their original single-line text.

# 3) The header cells now need to wrap their multi-line text (G1 stays
#    single line / unwrapped, matching the corporate-name column which
#    never had an embedded example).
$ws.Range("A1:F1").WrapText = $true
$ws.Range("H1:M1").WrapText = $true

# 4) Grow the header row to fit the now much taller wrapped text.
$ws.Rows.Item(1).RowHeight = 313.2

# 5) Swap the sample applicant name.
$ws.Range("B2").Value = "홍길동"

# 6) Widen the "교육이수번호/인증서번호" column now that its header text
#    is longer.
$ws.Columns.Item(9).ColumnWidth = 22.14

# 7) Restore the (somewhat arbitrary) last-saved selection.
$ws.Range("F6").Select()
